$wb = $excel.ActiveWorkbook

# --- Sheet: funding --- (text content itself doesn't change, only the
#     shared-string table shrinks because "HR&L" is removed elsewhere;
#     COM will handle the shared-string table automatically once we
#     make the corresponding edits below.)

# --- Sheet: keyword_set ---
$ws = $wb.Worksheets.Item("keyword_set")

# Row 8 used to hold "HR&L"; that duplicate keyword is removed and the
# existing "Healthy Rivers and Landscapes" keyword now sits at row 8.
$ws.Range("A8").Value = "Healthy Rivers and Landscapes"

# New fish-species keywords appended (row 9 first so shared-string order
# matches the target workbook).
$ws.Range("A9").Value = "o. mykiss"

# Row 2: "chinook" -> "chinook salmon"
$ws.Range("A2").Value = "chinook salmon"

$ws.Range("A10").Value = "pikeminnow"
$ws.Range("A11").Value = "sacramento sucker "
$ws.Range("A12").Value = "cypriniform"
$ws.Range("A13").Value = "tule perch"

$ws.Range("C18").Select()
